$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("18026822-a885-43a8-884a-3682ab6bf9bd", "Marcos", "joaovitorbatista12337@gmail.com"),
    @("cb783894-e9a6-4eee-887a-d9dbeca80720", "Joaquim", "email.com"),
    @("2b800868-aa5c-4131-b8c8-d0fd4e56494d", "VemPAo", "Pao@gmail.com")
)

$row = 3
foreach ($entry in $data) {
    $ws.Cells.Item($row, 1).Value = $entry[0]
    $ws.Cells.Item($row, 2).Value = $entry[1]
    $ws.Cells.Item($row, 3).Value = $entry[2]
    $row++
}
